$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: remove the long explanatory note from A1, unmerge A1:N1, ---
# --- and reset its alignment to General (creates the new "applyAlignment, no horizontal" style) ---
$ws.Range("A1:N1").UnMerge()
$ws.Range("A1").Value = ""
$ws.Range("A1:N1").HorizontalAlignment = 1

# --- New DS2 Anemometer section (rows 41-44), mirroring the existing anemometer/soil sections ---
# (cell values are written in the same order the original author entered them, so the
#  shared-string table gets built up in the same order as the target workbook)
$ws.Range("A41").Value = "DS2 Anemometer:"
$ws.Range("A41").Font.Bold = $true
$ws.Range("B41").Value = "No itermediate wire necessary"

$ws.Range("C43").Value = """SONIC DATA"" on LHS of LEMS Shield"
$ws.Range("C44").Value = """5V"" on LHS of LEMS Shield"
$ws.Range("C42").Value = """GND"" on LHS of LEMS Shield"

$ws.Range("A42").Value = "Shield"
$ws.Range("B42").Value = "Screw"

$ws.Range("A43").Value = "Red"
$ws.Range("B43").Value = "Screw"

$ws.Range("A44").Value = "White"
$ws.Range("B44").Value = "Screw"

# --- "OR" label (A21) becomes italic instead of bold ---
$ws.Range("A21").Font.Bold = $false
$ws.Range("A21").Font.Italic = $true

# --- Column C needs to be wide enough for the new, longer text ---
$ws.Columns.Item(3).ColumnWidth = 30.14

# --- Update the view: scroll down a bit and select C43 ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C43").Select() | Out-Null
